$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.165483474731445
$ws.Range("B1").Value = 2.80886173248291
$ws.Range("C1").Value = 5.080277442932129
$ws.Range("D1").Value = 2.088796615600586
$ws.Range("E1").Value = 1.164832472801208
